# Update want-to-go counts (column F) across all four worksheets,
# matching the gh-pages data refresh (commit 456a3b4).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 1165
$ws.Range("F4").Value = 1233
$ws.Range("F6").Value = 167
$ws.Range("F7").Value = 536
$ws.Range("F8").Value = 314
$ws.Range("F10").Value = 1256
$ws.Range("F11").Value = 28471
$ws.Range("F12").Value = 3379
$ws.Range("F13").Value = 34
$ws.Range("F14").Value = 249
$ws.Range("F15").Value = 464
$ws.Range("F16").Value = 19
$ws.Range("F18").Value = 8
$ws.Range("F19").Value = 321
$ws.Range("F20").Value = 607
$ws.Range("F21").Value = 265
$ws.Range("F22").Value = 260
$ws.Range("F23").Value = 347
$ws.Range("F25").Value = 43
$ws.Range("F26").Value = 656
$ws.Range("F28").Value = 95
$ws.Range("F29").Value = 522
$ws.Range("F30").Value = 70
$ws.Range("F31").Value = 33
$ws.Range("F32").Value = 611
$ws.Range("F33").Value = 236
$ws.Range("F34").Value = 40

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F6").Value = 378
$ws.Range("F7").Value = 820
$ws.Range("F9").Value = 87
$ws.Range("F10").Value = 268
$ws.Range("F11").Value = 4240
$ws.Range("F13").Value = 185
$ws.Range("F18").Value = 34
$ws.Range("F22").Value = 4238

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 250
$ws.Range("F4").Value = 1163
$ws.Range("F5").Value = 292

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 250
$ws.Range("F4").Value = 1163
$ws.Range("F7").Value = 378
$ws.Range("F8").Value = 292
$ws.Range("F9").Value = 820
$ws.Range("F10").Value = 1165
$ws.Range("F11").Value = 1233
$ws.Range("F12").Value = 167
$ws.Range("F13").Value = 536
$ws.Range("F14").Value = 314
$ws.Range("F17").Value = 1256
$ws.Range("F18").Value = 87
$ws.Range("F19").Value = 87
$ws.Range("F20").Value = 268
$ws.Range("F22").Value = 185
$ws.Range("F25").Value = 3379
$ws.Range("F26").Value = 249
$ws.Range("F29").Value = 464
$ws.Range("F30").Value = 19
$ws.Range("F32").Value = 34
$ws.Range("F33").Value = 321
$ws.Range("F34").Value = 607
$ws.Range("F35").Value = 265
$ws.Range("F36").Value = 347
$ws.Range("F38").Value = 44
$ws.Range("F39").Value = 656
$ws.Range("F42").Value = 95
$ws.Range("F45").Value = 70
$ws.Range("F46").Value = 33
$ws.Range("F47").Value = 611
$ws.Range("F48").Value = 236
$ws.Range("F49").Value = 40

